# Add a new "CatchShare" (TURF upside) column to the data table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column D
$ws.Range("D1").Value = "CatchShare"

# Fill D2:D145 with 1 (numeric flag, matching every existing data row)
$ws.Range("D2:D145").Value = 1

# Best-effort restore of window/view state seen in the authored edit.
$win = $excel.ActiveWindow
$win.Left = 1620
$win.Top = 100

$ws.Range("A138").Select()
$ws.Range("D152").Select()
